$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row's bold style is only used by the header; clearing the
# bold flag first merges its formatting into the (already-present)
# non-bold data-row style before the header row disappears.
$ws.Range("A1:C1").Font.Bold = $false

# Drop the header row ("Items"/"Price"/"Discount") - the data shifts up
# to start at row 1.
$ws.Rows.Item(1).Delete()

# Drop the "Discount" column - only Items/Price remain.
$ws.Columns.Item(3).Delete()

# New item names and prices (two extra rows: Carrots, Cake).
$items  = @("Gum", "Chips", "Bread", "Hot Chips", "Carrots", "Cake")
$prices = @(5.0, 1.0, 8.0, 9.0, 2.0, 3.0)

for ($i = 0; $i -lt $items.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $items[$i]
    $ws.Cells.Item($row, 2).Value = $prices[$i]
}

# Rows 5-6 are brand new (beyond the original 4-row data range) so they
# don't inherit any formatting automatically - copy the style from the
# last pre-existing data row (row 4) onto them.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
